$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5 (Personal Income Tax) re-executed and still Pass/Y, with new timestamps
$ws.Range("B2").Value = "Thu Nov 07 18:00:26 EST 2024"
$ws.Range("B3").Value = "Thu Nov 07 18:00:40 EST 2024"
$ws.Range("B4").Value = "Thu Nov 07 18:00:54 EST 2024"
$ws.Range("B5").Value = "Thu Nov 07 18:01:09 EST 2024"

# Rows 6-7 (Estate Tax) removed from execution: Result=Fail, Execute=DoNotRun, new timestamps
$ws.Range("A6").Value = "Fail"
$ws.Range("B6").Value = "Thu Nov 07 16:44:02 EST 2024"
$ws.Range("C6").Value = "DoNotRun"

$ws.Range("A7").Value = "Fail"
$ws.Range("B7").Value = "Thu Nov 07 16:44:17 EST 2024"
$ws.Range("C7").Value = "DoNotRun"

# Update selection to match the new active cell/selection range
$ws.Range("C6:C7").Select() | Out-Null
